# Mise à jour de certains champs de Modules et de Professeurs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header: "Matières enseignés" (adds a new shared string,
# extends dimension/row span to E1 automatically)
$ws.Range("E1").Value = "Matières enseignés"

# Column widths for C, D, E (values chosen so the stored OOXML width,
# after this runtime's ColumnWidth -> width rounding, lands as close as
# possible to the authored widths 27.5703125 / 15.7109375 / 31.7109375)
$ws.Columns.Item(3).ColumnWidth = 26.6666666666667
$ws.Columns.Item(4).ColumnWidth = 14.8333333333333
$ws.Columns.Item(5).ColumnWidth = 30.8333333333333

# Move/collapse the selection to E6, matching the saved sheet view state
$ws.Range("E6").Select()
